# Add a new "CRYP" / "CryptoTransaction" entry to the ExternalPurpose1Code
# code set, as the new last row of that group (immediately after the
# existing "DEBT" / "ChargesBorneByDebtor" row, and before the first
# "ExternalSystemBalanceType1Code" row). This pushes every following row
# down by one and grows the table/used range from E692 to E693.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 614, pushing rows 614..692 down to 615..693.
$ws.Rows("614:614").Insert()

# Populate the new row with the CRYP code set entry.
$ws.Range("A614").Value = "ExternalPurpose1Code"
$ws.Range("B614").Value = "CRYP"
$ws.Range("C614").Value = "CryptoTransaction"
$ws.Range("D614").Value = "Transaction is for the purchase of cryptocurrency"

# Match the wrap-text formatting used by the rest of column D, and the
# single-line row height (17) used by similar one-line entries.
$ws.Range("D614").WrapText = $true
$ws.Rows("614:614").RowHeight = 17

# The worksheet table (Table1) covered A1:E692; grow it to A1:E693 so the
# new row is included, matching the autoFilter/table ref update.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E693"))

# Reflect the updated selection/active cell seen in the saved view state.
$ws.Activate()
$ws.Range("D613").Select()
$excel.ActiveWindow.ScrollRow = 589
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save()
